$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.109.09"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +4.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.899.56"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.63%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.97"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4998"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.39%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.85"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2979"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +7.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06689"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.898.83"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "17.11"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07269"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6830"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +6.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "86.41"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.90%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.901"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +5.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.070.92"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.99%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008041"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +10.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.98"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +6.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.148.55"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +5.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.0000"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.792"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +5.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.708"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +6.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.311"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +6.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "148.79"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "132.80"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.89"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.975"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +5.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.387"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.248"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08794"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +5.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.961"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05120"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +4.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.132"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7075"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.688"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.792"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.238"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9572"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01674"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +5.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.020"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.78%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4270"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.71%  "
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "103.70"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.518"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1269"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.89%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05758"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.09"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.75%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.334"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.00%  "
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3773"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +4.57%  "
